$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.1526236666666667
$ws.Range("H2").Value = 0.457871
$ws.Range("M2").Value = 27.30988266666667
$ws.Range("N2").Value = 81.929648
$ws.Range("O2").Value = 0.2168690090390243
$ws.Range("P2").Value = 0.2168690090390242
$ws.Range("Q2").Value = 4.168134428823111
$ws.Range("R2").Value = 37.513209859408
$ws.Range("S2").Value = 0.2168690090390243
$ws.Range("T2").Value = 0.2168690090390242

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1526236666666667
$ws.Range("H3").Value = 0.457871
$ws.Range("N3").Value = 61.114908
$ws.Range("O3").Value = 0.1617720795708915
$ws.Range("P3").Value = 0.1617720795708915
$ws.Range("Q3").Value = 3.109193782318666
$ws.Range("R3").Value = 27.982744040868
$ws.Range("S3").Value = 0.1617720795708915
$ws.Range("T3").Value = 0.1617720795708915

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1526236666666667
$ws.Range("H4").Value = 0.457871
$ws.Range("M4").Value = 42.52135866666666
$ws.Range("N4").Value = 127.564076
$ws.Range("O4").Value = 0.3376640254953711
$ws.Range("P4").Value = 0.3376640254953711
$ws.Range("Q4").Value = 6.48976567135511
$ws.Range("R4").Value = 58.407891042196
$ws.Range("S4").Value = 0.3376640254953711
$ws.Range("T4").Value = 0.3376640254953711

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1526236666666667
$ws.Range("H5").Value = 0.457871
$ws.Range("M5").Value = 6.417394333333333
$ws.Range("N5").Value = 19.252183
$ws.Range("O5").Value = 0.05096081761571768
$ws.Range("P5").Value = 0.05096081761571768
$ws.Range("Q5").Value = 0.9794462535992221
$ws.Range("R5").Value = 8.815016282393
$ws.Range("S5").Value = 0.05096081761571768
$ws.Range("T5").Value = 0.05096081761571768

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.1526236666666667
$ws.Range("H6").Value = 0.457871
$ws.Range("M6").Value = 29.307738
$ws.Range("N6").Value = 87.923214
$ws.Range("O6").Value = 0.2327340682789955
$ws.Range("P6").Value = 0.2327340682789955
$ws.Range("Q6").Value = 4.473054435266
$ws.Range("R6").Value = 40.257489917394
$ws.Range("S6").Value = 0.2327340682789955
$ws.Range("T6").Value = 0.2327340682789955

$wb.Save()
